# Apply the May-23-2023 crypto price/volume refresh described in the commit.
# D-column price strings that are valid numeric literals (e.g. "1.001") would be
# auto-coerced to numbers by Excel on assignment, so those are written with a
# leading apostrophe (Excel's standard force-text prefix) to keep them as text,
# exactly like the source data (which stores every Price/Volume cell as a string).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '27.322.37'
$ws.Range("E2").Value = '  +1.17%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.855.80'
$ws.Range("E3").Value = '  +1.45%  '

# Row 4: TetherUSD
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.87%  '

# Row 5: BNB
$ws.Range("D5").Value = '''313.78'
$ws.Range("E5").Value = '  +0.63%  '

# Row 6: USDC
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  -0.74%  '

# Row 7: XRP
$ws.Range("D7").Value = '''0.4606'
$ws.Range("E7").Value = '  -0.90%  '

# Row 8: Cardano
$ws.Range("D8").Value = '''0.3710'
$ws.Range("E8").Value = '  +0.04%  '

# Row 9: Dogecoin
$ws.Range("D9").Value = '''0.07319'
$ws.Range("E9").Value = '  -0.67%  '

# Row 10: Polygon
$ws.Range("D10").Value = '''0.8798'

# Row 11: Solana
$ws.Range("E11").Value = '  -0.63%  '

# Row 12: TRON
$ws.Range("D12").Value = '''0.07804'
$ws.Range("E12").Value = '  -0.89%  '

# Row 13: WrappedEther
$ws.Range("D13").Value = '1.813.12'
$ws.Range("E13").Value = '  +2.21%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '''5.388'
$ws.Range("E14").Value = '  +0.51%  '

# Row 15: Chainlink
$ws.Range("D15").Value = '''6.541'

# Row 16: Litecoin
$ws.Range("D16").Value = '''91.81'

# Row 17: BinanceUSD
$ws.Range("D17").Value = '''1.002'
$ws.Range("E17").Value = '  -0.80%  '

# Row 18: ShibaInu
$ws.Range("D18").Value = '''0.000009033'
$ws.Range("E18").Value = '  +1.62%  '

# Row 19: Dai
$ws.Range("D19").Value = '''1.000'
$ws.Range("E19").Value = '  -0.72%  '

# Row 20: Avalanche
$ws.Range("E20").Value = '  +0.38%  '

# Row 21: WrappedBTC
$ws.Range("D21").Value = '27.347.52'
$ws.Range("E21").Value = '  +1.37%  '

# Row 22: Uniswap
$ws.Range("D22").Value = '''5.129'
$ws.Range("E22").Value = '  -0.70%  '

# Row 23: Cosmos
$ws.Range("D23").Value = '''10.53'
$ws.Range("E23").Value = '  -0.56%  '

# Row 24: WrappedliquidstakedEther2.0
$ws.Range("D24").Value = '2.128.35'
$ws.Range("E24").Value = '  +7.20%  '

# Row 25: Toncoin
$ws.Range("D25").Value = '''1.919'
$ws.Range("E25").Value = '  +4.64%  '

# Row 26: Monero
$ws.Range("D26").Value = '''152.10'
$ws.Range("E26").Value = '  -0.34%  '

# Row 27: EthereumClassic
$ws.Range("D27").Value = '''18.38'
$ws.Range("E27").Value = '  +0.50%  '

# Row 28: LidoDAOToken
$ws.Range("D28").Value = '''2.070'
$ws.Range("E28").Value = '  -1.45%  '

# Row 29: InternetComputer(DFINITY)
$ws.Range("D29").Value = '''5.106'
$ws.Range("E29").Value = '  -0.48%  '

# Row 30: BitcoinCash
$ws.Range("D30").Value = '''116.19'
$ws.Range("E30").Value = '  +0.50%  '

# Row 31: Stellar
$ws.Range("D31").Value = '''0.08861'
$ws.Range("E31").Value = '  -0.25%  '

# Row 32: ImmutableX
$ws.Range("D32").Value = '''0.7741'
$ws.Range("E32").Value = '  +6.30%  '

# Row 33: HuobiToken
$ws.Range("D33").Value = '''3.031'
$ws.Range("E33").Value = '  +1.58%  '

# Row 34: ARBITRUM
$ws.Range("D34").Value = '''1.177'
$ws.Range("E34").Value = '  +3.59%  '

# Row 35: Filecoin
$ws.Range("D35").Value = '''4.488'
$ws.Range("E35").Value = '  +0.82%  '

# Row 36: RenderToken
$ws.Range("D36").Value = '''2.633'
$ws.Range("E36").Value = '  +5.75%  '

# Row 37: VeChain
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '''1.078'
$ws.Range("E37").Value = '  -0.23%  '

# Row 38: TrustWalletToken
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.01960'
$ws.Range("E38").Value = '  +0.21%  '

# Row 39: Hedera
$ws.Range("D39").Value = '''0.05231'
$ws.Range("E39").Value = '  -0.33%  '

# Row 40: MXToken
$ws.Range("D40").Value = '''2.952'
$ws.Range("E40").Value = '  +0.68%  '

# Row 41: FraxShare
$ws.Range("D41").Value = '''7.049'
$ws.Range("E41").Value = '  -3.75%  '

# Row 42: TheSandbox
$ws.Range("D42").Value = '''0.5142'
$ws.Range("E42").Value = '  -1.13%  '

# Row 43: Algorand
$ws.Range("D43").Value = '''0.1637'
$ws.Range("E43").Value = '  +0.59%  '

# Row 44: Aptos
$ws.Range("D44").Value = '''8.390'
$ws.Range("E44").Value = '  +2.01%  '

# Row 45: Decentraland
$ws.Range("D45").Value = '''0.4832'
$ws.Range("E45").Value = '  -0.34%  '

# Row 46: EnergySwap
$ws.Range("E46").Value = '  +0.52%  '

# Row 47: PaxDollar
$ws.Range("D47").Value = '''1.000'
$ws.Range("E47").Value = '  -0.81%  '

# Row 48: Quant
$ws.Range("D48").Value = '''103.23'
$ws.Range("E48").Value = '  +0.44%  '

# Row 49: NEARProtocol
$ws.Range("D49").Value = '''1.652'
$ws.Range("E49").Value = '  +1.57%  '

# Row 50: Cronos
$ws.Range("D50").Value = '''0.06218'
$ws.Range("E50").Value = '  -0.33%  '

# Row 51: Aave
$ws.Range("D51").Value = '''65.86'
$ws.Range("E51").Value = '  +2.21%  '
